# Natmi following Dr Hou advice
# Sirpa-Cd47 sheet: the NATMI pipeline was rerun with an updated minimum
# "expressing cells" threshold (1 -> 3), which changes the downstream
# expression / specificity / edge-weight statistics for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updated values, keyed by column letter -> new value.
$rowUpdates = @{
    2  = @{ E=3; G=94.95332266666666;  H=284.859968;          I=0.6912729447872573;  J=0.6912729447872574;
            K=3; M=59.45197733333333;  N=178.355932;           O=0.304222453049858;   P=0.304222453049858;
            Q=5645.162786903353;       R=50806.46508213017;    S=0.2103007509901784;  T=0.2103007509901785 }
    3  = @{ E=3; G=94.95332266666666;  H=284.859968;          I=0.6912729447872573;  J=0.6912729447872574;
            K=3; M=107.1770123333333;  N=321.531037;           O=0.548436823552382;   P=0.5484368235523819;
            Q=10176.81343453631;       R=91591.32091082681;    S=0.3791195380468245;  T=0.3791195380468245 }
    4  = @{ E=3; G=94.95332266666666;  H=284.859968;          I=0.6912729447872573;  J=0.6912729447872574;
            K=3; M=28.793724;          N=86.38117199999999;    O=0.1473407233977601;  P=0.1473407233977601;
            Q=2734.059765746944;       R=24606.53789172249;    S=0.1018526557502544;  T=0.1018526557502544 }
    5  = @{ E=3; G=36.59611566666667;  H=109.788347;           I=0.2664246382770613;  J=0.2664246382770613;
            K=3; M=59.45197733333333;  N=178.355932;           O=0.304222453049858;   P=0.304222453049858;
            Q=2175.711439102712;       R=19581.4029519244;     S=0.08105235700956866; T=0.08105235700956867 }
    6  = @{ E=3; G=36.59611566666667;  H=109.788347;           I=0.2664246382770613;  J=0.2664246382770613;
            K=3; M=107.1770123333333;  N=321.531037;           O=0.548436823552382;   P=0.5484368235523819;
            Q=3922.262340158426;       R=35300.36106142584;    S=0.1461170823327639;  T=0.1461170823327639 }
    7  = @{ E=3; G=36.59611566666667;  H=109.788347;           I=0.2664246382770613;  J=0.2664246382770613;
            K=3; M=28.793724;          N=86.38117199999999;    O=0.1473407233977601;  P=0.1473407233977601;
            Q=1053.738453978076;       R=9483.646085802684;    S=0.03925519893472876; T=0.03925519893472876 }
    8  = @{ E=3; G=5.810664333333333;  H=17.431993;             I=0.04230241693568138; J=0.04230241693568138;
            K=3; M=59.45197733333333;  N=178.355932;           O=0.304222453049858;   P=0.304222453049858;
            Q=345.4554842369417;       R=3109.099358132476;    S=0.01286934505011084; T=0.01286934505011085 }
    9  = @{ E=3; G=5.810664333333333;  H=17.431993;             I=0.04230241693568138; J=0.04230241693568138;
            K=3; M=107.1770123333333;  N=321.531037;           O=0.548436823552382;   P=0.5484368235523819;
            Q=622.7696429185266;       R=5604.92678626674;     S=0.02320020317279358; T=0.02320020317279358 }
    10 = @{ E=3; G=5.810664333333333;  H=17.431993;             I=0.04230241693568138; J=0.04230241693568138;
            K=3; M=28.793724;          N=86.38117199999999;    O=0.1473407233977601;  P=0.1473407233977601;
            Q=167.3106650706439;       R=1505.795985635796;    S=0.006232868712776951; T=0.006232868712776951 }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
